$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute("Systems analyst,")
$start = $r.Start
$tRng = $d.Range($start + 14, $start + 15)
$tRng.Text = "is"
$isRng = $d.Range($start + 14, $start + 16)
$isRng.Italic = 1
Write-Output "xml1=$($isRng.WordOpenXML)"
$isRng.Italic = 0
Write-Output "xml2=$($isRng.WordOpenXML)"
Write-Output "final=[$($d.Content.Text)]"
